$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039448040340348
$ws.Range("D2").Value = 1.047968690964252
$ws.Range("E2").Value = 1.037949812231992
$ws.Range("F2").Value = 1.05267353513313
$ws.Range("I2").Value = 1.026257537234071
$ws.Range("J2").Value = 1.044540002810251
$ws.Range("K2").Value = 1.050729846239048
$ws.Range("L2").Value = 1.040739236913732
$ws.Range("M2").Value = 1.055421616441516
$ws.Range("N2").Value = 1.018665458093369
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041921567296975
$ws.Range("D3").Value = 1.050447336200944
$ws.Range("E3").Value = 1.040086030807114
$ws.Range("F3").Value = 1.055370732492939
$ws.Range("I3").Value = 1.026539659425466
$ws.Range("J3").Value = 1.046651127076961
$ws.Range("K3").Value = 1.053016131983251
$ws.Range("L3").Value = 1.04268186526578
$ws.Range("M3").Value = 1.057926880983287
$ws.Range("N3").Value = 1.019366059093507
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043503413442479
$ws.Range("D4").Value = 1.052029851983047
$ws.Range("E4").Value = 1.041452462985078
$ws.Range("F4").Value = 1.057086737650216
$ws.Range("I4").Value = 1.026711049692088
$ws.Range("J4").Value = 1.047998618626887
$ws.Range("K4").Value = 1.054473943470437
$ws.Range("L4").Value = 1.043922716715861
$ws.Range("M4").Value = 1.059518520267879
$ws.Range("N4").Value = 1.019813134988916
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044164043265282
$ws.Range("D5").Value = 1.052690136028571
$ws.Range("E5").Value = 1.04202319931339
$ws.Range("F5").Value = 1.057801264997801
$ws.Range("I5").Value = 1.026780449448248
$ws.Range("J5").Value = 1.048560751106046
$ws.Range("K5").Value = 1.055081741414345
$ws.Range("L5").Value = 1.04444058030021
$ws.Range("M5").Value = 1.060180717635613
$ws.Range("N5").Value = 1.019999616154784
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044274711903928
$ws.Range("D6").Value = 1.052800709969924
$ws.Range("E6").Value = 1.04211881309755
$ws.Range("F6").Value = 1.057920837268245
$ws.Range("I6").Value = 1.026791946953355
$ws.Range("J6").Value = 1.048654882829174
$ws.Range("K6").Value = 1.055183499045111
$ws.Range("L6").Value = 1.044527311816191
$ws.Range("M6").Value = 1.060291500590408
$ws.Range("N6").Value = 1.020030841810864
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043512257888361
$ws.Range("D7").Value = 1.052038694266727
$ws.Range("E7").Value = 1.04146010367032
$ws.Range("F7").Value = 1.057096312076465
$ws.Range("I7").Value = 1.026711987413249
$ws.Range("J7").Value = 1.04800614685901
$ws.Range("K7").Value = 1.054482084675173
$ws.Range("L7").Value = 1.043929651226026
$ws.Range("M7").Value = 1.059527395648269
$ws.Range("N7").Value = 1.019815632495742
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040287909814005
$ws.Range("D8").Value = 1.04881084058223
$ws.Range("E8").Value = 1.038675087987094
$ws.Range("F8").Value = 1.053591201055903
$ws.Range("I8").Value = 1.02635520499912
$ws.Range("J8").Value = 1.045257361145016
$ws.Range("K8").Value = 1.051507033136384
$ws.Range("L8").Value = 1.041399151192251
$ws.Range("M8").Value = 1.056274450277394
$ws.Range("N8").Value = 1.018903543511675
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034458514301874
$ws.Range("D9").Value = 1.042954772114192
$ws.Range("E9").Value = 1.033642363575592
$ws.Range("F9").Value = 1.04718476843715
$ws.Range("I9").Value = 1.02564007415476
$ws.Range("J9").Value = 1.040267493734139
$ws.Range("K9").Value = 1.046094828200115
$ws.Range("L9").Value = 1.036812651496181
$ws.Range("M9").Value = 1.05031128661216
$ws.Range("N9").Value = 1.017247015695852
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030466323579439
$ws.Range("D10").Value = 1.038930622308856
$ws.Range("E10").Value = 1.030197452780379
$ws.Range("F10").Value = 1.042750543817115
$ws.Range("I10").Value = 1.025103760967125
$ws.Range("J10").Value = 1.036836619343741
$ws.Range("K10").Value = 1.042365782689046
$ws.Range("L10").Value = 1.033663912088359
$ws.Range("M10").Value = 1.046172195616202
$ws.Range("N10").Value = 1.01610750393481
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028710969747995
$ws.Range("D11").Value = 1.037157968646538
$ws.Range("E11").Value = 1.028683154046165
$ws.Range("F11").Value = 1.04078969444646
$ws.Range("I11").Value = 1.024857048342973
$ws.Range("J11").Value = 1.035324847514745
$ws.Range("K11").Value = 1.040720784721045
$ws.Range("L11").Value = 1.032277596009187
$ws.Range("M11").Value = 1.04433913215651
$ws.Range("N11").Value = 1.015605268625014
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028054799442329
$ws.Range("D12").Value = 1.036494843193883
$ws.Range("E12").Value = 1.028117156975694
$ws.Range("F12").Value = 1.040055031816155
$ws.Range("I12").Value = 1.024763198898323
$ws.Range("J12").Value = 1.034759246041531
$ws.Range("K12").Value = 1.040105062843543
$ws.Range("L12").Value = 1.031759102354479
$ws.Range("M12").Value = 1.043651940593599
$ws.Range("N12").Value = 1.01541734829441
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028195740505922
$ws.Range("D13").Value = 1.036637300263217
$ws.Range("E13").Value = 1.028238726494519
$ws.Range("F13").Value = 1.040212908343927
$ws.Range("I13").Value = 1.024783430458344
$ws.Range("J13").Value = 1.034880755434956
$ws.Range("K13").Value = 1.04023735220156
$ws.Range("L13").Value = 1.031870483745205
$ws.Range("M13").Value = 1.043799634043055
$ws.Range("N13").Value = 1.015457720448397
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028656816081841
$ws.Range("D14").Value = 1.037103250917281
$ws.Range("E14").Value = 1.028636441051582
$ws.Range("F14").Value = 1.040729096957642
$ws.Range("I14").Value = 1.024849335989125
$ws.Range("J14").Value = 1.035278178354792
$ws.Range("K14").Value = 1.040669985627275
$ws.Range("L14").Value = 1.032234810351492
$ws.Range("M14").Value = 1.044282458483255
$ws.Range("N14").Value = 1.015589763229377
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028940345215785
$ws.Range("D15").Value = 1.037389713334808
$ws.Range("E15").Value = 1.028881016135274
$ws.Range("F15").Value = 1.041046295156189
$ws.Range("I15").Value = 1.024889648767836
$ws.Range("J15").Value = 1.035522501502542
$ws.Range("K15").Value = 1.040935918590296
$ws.Range("L15").Value = 1.032458809503043
$ws.Range("M15").Value = 1.044579100881955
$ws.Range("N15").Value = 1.015670936585203
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030582246023409
$ws.Range("D16").Value = 1.039047619189421
$ws.Range("E16").Value = 1.030297464972587
$ws.Range("F16").Value = 1.042879803646663
$ws.Range("I16").Value = 1.025119826502772
$ws.Range("J16").Value = 1.036936388220053
$ws.Range("K16").Value = 1.042474305288322
$ws.Range("L16").Value = 1.033755425342913
$ws.Range("M16").Value = 1.046292975122265
$ws.Range("N16").Value = 1.016140646223215
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031604918855937
$ws.Range("D17").Value = 1.040079398581542
$ws.Range("E17").Value = 1.031179824933258
$ws.Range("F17").Value = 1.044018863225097
$ws.Range("I17").Value = 1.025260310053453
$ws.Range("J17").Value = 1.037816183804598
$ws.Range("K17").Value = 1.043431083491169
$ws.Range("L17").Value = 1.034562550171299
$ws.Range("M17").Value = 1.047356993795896
$ws.Range("N17").Value = 1.016432891967734
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032198859475655
$ws.Range("D18").Value = 1.040678317590979
$ws.Range("E18").Value = 1.031692315496394
$ws.Range("F18").Value = 1.044679333784131
$ws.Range("I18").Value = 1.025340856449817
$ws.Range("J18").Value = 1.038326836669568
$ws.Range("K18").Value = 1.043986242993921
$ws.Range("L18").Value = 1.035031131828732
$ws.Range("M18").Value = 1.047973692919201
$ws.Range("N18").Value = 1.016602506041543
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032400946616437
$ws.Range("D19").Value = 1.040882045966839
$ws.Range("E19").Value = 1.031866696129803
$ws.Range("F19").Value = 1.044903877152071
$ws.Range("I19").Value = 1.025368084994579
$ws.Range("J19").Value = 1.038500532975645
$ws.Range("K19").Value = 1.044175048209141
$ws.Range("L19").Value = 1.035190536203741
$ws.Range("M19").Value = 1.048183310921318
$ws.Range("N19").Value = 1.016660197499441
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031495462269007
$ws.Range("D20").Value = 1.039968999487407
$ws.Range("E20").Value = 1.031085381881976
$ws.Range("F20").Value = 1.043897060065234
$ws.Range("I20").Value = 1.02524538205107
$ws.Range("J20").Value = 1.037722051322105
$ws.Range("K20").Value = 1.043328732576664
$ws.Range("L20").Value = 1.034476181747686
$ws.Range("M20").Value = 1.047243241904128
$ws.Range("N20").Value = 1.016401624772669
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028521156646004
$ws.Range("D21").Value = 1.036966170578511
$ws.Range("E21").Value = 1.028519422157594
$ws.Range("F21").Value = 1.040577268166118
$ws.Range("I21").Value = 1.02482998971969
$ws.Range("J21").Value = 1.035161260409669
$ws.Range("K21").Value = 1.040542716707091
$ws.Range("L21").Value = 1.032127624309313
$ws.Range("M21").Value = 1.044140454394482
$ws.Range("N21").Value = 1.01555091803264
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026626985586991
$ws.Range("D22").Value = 1.035050998865835
$ws.Range("E22").Value = 1.026885677726828
$ws.Range("F22").Value = 1.03845334242485
$ws.Range("I22").Value = 1.02455601521493
$ws.Range("J22").Value = 1.033527618997447
$ws.Range("K22").Value = 1.038763789129509
$ws.Range("L22").Value = 1.030630368798218
$ws.Range("M22").Value = 1.042153004517127
$ws.Range("N22").Value = 1.015008108688124
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027633455457342
$ws.Range("D23").Value = 1.036068895531224
$ws.Range("E23").Value = 1.027753733836509
$ws.Range("F23").Value = 1.039582813407761
$ws.Range("I23").Value = 1.024702479290827
$ws.Range("J23").Value = 1.034395922109455
$ws.Range("K23").Value = 1.03970946538703
$ws.Range("L23").Value = 1.031426087100594
$ws.Range("M23").Value = 1.043210121179999
$ws.Range("N23").Value = 1.015296629298061
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031544928912994
$ws.Range("D24").Value = 1.04001889303917
$ws.Range("E24").Value = 1.031128063345321
$ws.Range("F24").Value = 1.043952109782708
$ws.Range("I24").Value = 1.025252131682261
$ws.Range("J24").Value = 1.037764593508743
$ws.Range("K24").Value = 1.043374989541681
$ws.Range("L24").Value = 1.034515214713455
$ws.Range("M24").Value = 1.047294653604701
$ws.Range("N24").Value = 1.016415755689706
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035983716831591
$ws.Range("D25").Value = 1.044489337925587
$ws.Range("E25").Value = 1.034958836697648
$ws.Range("F25").Value = 1.048869099975823
$ws.Range("I25").Value = 1.025835337504291
$ws.Range("J25").Value = 1.041575416159012
$ws.Range("K25").Value = 1.046094828200115
$ws.Range("L25").Value = 1.038014011867038
$ws.Range("M25").Value = 1.051881115388859
$ws.Range("N25").Value = 1.017681311744702
